$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "Aufhol-Woche" topic row into the course schedule, shifting the
# remaining topics down by one (the last topic, "Fallstudien", drops off the
# list), and update the comment texts accordingly.

$ws.Range("B6").Value = "Aufhol-Woche"
$ws.Range("D6").Value = "Am Di., 1.11. entfällt die Vorlesung. Am Do., 3. 11. entfällt die Übung."

$ws.Range("B7").Value = "Die Post befragen"
$ws.Range("D7").Value = "Ab diese Woche benötigen wir rstanarm."

$ws.Range("B9").Value = "Gauss-Modelle"
$ws.Range("B10").Value = "Lineare Modelle"
$ws.Range("B11").Value = "Metrische AV"
$ws.Range("B12").Value = "Kausalinferenz 1"
$ws.Range("B13").Value = "Kausalinferenz 2"

$ws.Range("B7").Select()
